# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "69.124.36"
Set-TextCell $ws.Range("E2") "  +2.14%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.776.29"
Set-TextCell $ws.Range("E3") "  +0.15%  "

# Row 4
Set-TextCell $ws.Range("E4") "  -0.38%  "

# Row 5
Set-TextCell $ws.Range("D5") "625.23"
Set-TextCell $ws.Range("E5") "  +4.40%  "

# Row 6
Set-TextCell $ws.Range("D6") "165.97"
Set-TextCell $ws.Range("E6") "  +1.92%  "

# Row 7
Set-TextCell $ws.Range("D7") "3.773.66"
Set-TextCell $ws.Range("E7") "  +0.15%  "

# Row 8
Set-TextCell $ws.Range("E8") "  -0.09%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.520"
Set-TextCell $ws.Range("E9") "  +1.45%  "

# Row 10
Set-TextCell $ws.Range("E10") "  +3.03%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.456"
Set-TextCell $ws.Range("E11") "  +2.85%  "

# Row 12
Set-TextCell $ws.Range("D12") "6.71"
Set-TextCell $ws.Range("E12") "  +1.89%  "

# Row 13
Set-TextCell $ws.Range("D13") "0.0000246"
Set-TextCell $ws.Range("E13") "  +1.06%  "

# Row 14
Set-TextCell $ws.Range("E14") "  +1.84%  "

# Row 15
Set-TextCell $ws.Range("D15") "4.414.47"
Set-TextCell $ws.Range("E15") "  +0.26%  "

# Row 16
Set-TextCell $ws.Range("D16") "3.767.89"
Set-TextCell $ws.Range("E16") "  +0.40%  "

# Row 17
Set-TextCell $ws.Range("D17") "69.151.02"
Set-TextCell $ws.Range("E17") "  +2.08%  "

# Row 18
Set-TextCell $ws.Range("D18") "17.68"
Set-TextCell $ws.Range("E18") "  -2.79%  "

# Row 19
Set-TextCell $ws.Range("D19") "7.11"
Set-TextCell $ws.Range("E19") "  +1.88%  "

# Row 20
Set-TextCell $ws.Range("E20") "  -0.96%  "

# Row 21
Set-TextCell $ws.Range("D21") "468.11"
Set-TextCell $ws.Range("E21") "  +2.51%  "

# Row 22
Set-TextCell $ws.Range("D22") "9.62"
Set-TextCell $ws.Range("E22") "  +1.75%  "

# Row 23
Set-TextCell $ws.Range("D23") "0.706"
Set-TextCell $ws.Range("E23") "  +2.30%  "

# Row 24
Set-TextCell $ws.Range("D24") "0.0000147"
Set-TextCell $ws.Range("E24") "  +3.81%  "

# Row 25
Set-TextCell $ws.Range("D25") "83.12"
Set-TextCell $ws.Range("E25") "  +0.56%  "

# Row 26
Set-TextCell $ws.Range("B26") "Fetch.AI"
Set-TextCell $ws.Range("C26") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Range("D26") "2.17"
Set-TextCell $ws.Range("E26") "  +4.10%  "

# Row 27
Set-TextCell $ws.Range("B27") "InternetComputer(DFINITY)"
Set-TextCell $ws.Range("C27") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D27") "12.01"
Set-TextCell $ws.Range("E27") "  +1.69%  "

# Row 28
Set-TextCell $ws.Range("B28") "RenderToken"
Set-TextCell $ws.Range("C28") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D28") "10.02"
Set-TextCell $ws.Range("E28") "  +2.16%  "

# Row 29
Set-TextCell $ws.Range("B29") "Dai"
Set-TextCell $ws.Range("C29") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws.Range("D29") "1.00"
Set-TextCell $ws.Range("E29") "  -0.10%  "

# Row 30
Set-TextCell $ws.Range("D30") "3.927.34"
Set-TextCell $ws.Range("E30") "  +0.25%  "

# Row 31
Set-TextCell $ws.Range("E31") "  +3.66%  "

# Row 32
Set-TextCell $ws.Range("D32") "2.24"
Set-TextCell $ws.Range("E32") "  +2.03%  "

# Row 33
Set-TextCell $ws.Range("D33") "7.25"

# Row 34
Set-TextCell $ws.Range("E34") "  -0.13%  "

# Row 35
Set-TextCell $ws.Range("E35") "  +0.13%  "

# Row 36
Set-TextCell $ws.Range("E36") "  +16.01%  "

# Row 37
Set-TextCell $ws.Range("B37") "Aptos"
Set-TextCell $ws.Range("C37") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws.Range("D37") "9.00"
Set-TextCell $ws.Range("E37") "  +0.87%  "

# Row 38
Set-TextCell $ws.Range("B38") "RenzoRestakedETH"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextCell $ws.Range("D38") "3.729.10"
Set-TextCell $ws.Range("E38") "  +0.26%  "

# Row 39
Set-TextCell $ws.Range("E39") "  +3.00%  "

# Row 40
Set-TextCell $ws.Range("E40") "  +8.24%  "

# Row 41
Set-TextCell $ws.Range("D41") "5.80"
Set-TextCell $ws.Range("E41") "  +0.98%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.967"
Set-TextCell $ws.Range("E42") "  -0.81%  "

# Row 43
Set-TextCell $ws.Range("E43") "  -0.10%  "

# Row 45
Set-TextCell $ws.Range("E45") "  +1.01%  "

# Row 46
Set-TextCell $ws.Range("D46") "43.20"
Set-TextCell $ws.Range("E46") "  -0.49%  "

# Row 47
Set-TextCell $ws.Range("E47") "  +4.77%  "

# Row 48
Set-TextCell $ws.Range("B48") "Monero"
Set-TextCell $ws.Range("C48") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D48") "151.95"
Set-TextCell $ws.Range("E48") "  -0.12%  "

# Row 49
Set-TextCell $ws.Range("B49") "OKB"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D49") "46.68"
Set-TextCell $ws.Range("E49") "  -0.91%  "

# Row 50
Set-TextCell $ws.Range("D50") "8.41"
Set-TextCell $ws.Range("E50") "  +1.82%  "

# Row 51
Set-TextCell $ws.Range("D51") "1.36"
Set-TextCell $ws.Range("E51") "  +0.40%  "
